# Add two new people (rows) to the "Tableau2" assignments table on Feuil1:
#   - DAVID / Laetitia     -> N/A / N/A        / ACT3 = 1      (inserted before DELPIERRE)
#   - MAJAJD / Ines        -> 1,2,3 / 1,2,3,4,5 / ACT3 = 1.2    (inserted before HALLEBARD)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert two blank worksheet rows at the right spots (before the old row 103
# "DELPIERRE" and before the old row 104 "HALLEBARD", which becomes row 106
# after the first insert) so everything below shifts down.
[void]$ws.Rows.Item(103).Insert()
[void]$ws.Rows.Item(105).Insert()

# The table range doesn't auto-grow when rows are spliced in directly, so
# resize it explicitly to cover the two new rows (A1:E114 -> A1:E116).
[void]$lo.Resize($ws.Range("A1:E116"))

# Row 103: DAVID / Laetitia. Fill PRENOM (B) before NOM (A) to match entry order.
$ws.Range("B103").Value2 = "Laetitia"
$ws.Range("A103").Value2 = "DAVID"
$ws.Range("C103").Value2 = "N/A"
$ws.Range("D103").Value2 = "N/A"
$ws.Range("E103").Value2 = 1
[void]$ws.Range("E103").Copy()
[void]$ws.Range("B103").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 105: MAJAJD / Ines.
$ws.Range("A105").Value2 = "MAJAJD"
$ws.Range("B105").Value2 = "Ines"
$ws.Range("C105").Value2 = "1,2,3"
$ws.Range("D105").Value2 = "1,2,3,4,5"
$ws.Range("E105").Value2 = 1.2

# Restore the selection to what it was left at in the saved file.
[void]$ws.Range("E106").Select()
